{"js": "// Update the two-digit x two-digit multiplication answers in the table.\n// Each old value is unique in the document, so a direct search + replace\n// per pair (processed in document order) is safe even though a couple of\n// the new values happen to collide with old values used elsewhere.\nconst pairs = [\n  [\"99\u00d722=2178\", \"97\u00d758=5626\"],\n  [\"63\u00d754=3402\", \"33\u00d722=726\"],\n  [\"21\u00d757=1197\", \"88\u00d753=4664\"],\n  [\"26\u00d771=1846\", \"98\u00d796=9408\"],\n  [\"39\u00d785=3315\", \"13\u00d789=1157\"],\n  [\"48\u00d723=1104\", \"61\u00d713=793\"],\n  [\"52\u00d746=2392\", \"14\u00d786=1204\"],\n  [\"71\u00d777=5467\", \"82\u00d791=7462\"],\n  [\"28\u00d727=756\", \"69\u00d734=2346\"],\n  [\"65\u00d753=3445\", \"19\u00d786=1634\"],\n  [\"38\u00d747=1786\", \"86\u00d792=7912\"],\n  [\"73\u00d774=5402\", \"81\u00d773=5913\"],\n  [\"32\u00d769=2208\", \"32\u00d759=1888\"],\n  [\"36\u00d754=1944\", \"73\u00d756=4088\"],\n  [\"15\u00d766=990\", \"98\u00d796=9408\"],\n  [\"43\u00d731=1333\", \"38\u00d792=3496\"],\n  [\"15\u00d738=570\", \"32\u00d769=2208\"],\n  [\"87\u00d769=6003\", \"26\u00d732=832\"],\n  [\"91\u00d772=6552\", \"87\u00d778=6786\"],\n  [\"14\u00d716=224\", \"36\u00d726=936\"],\n  [\"75\u00d776=5700\", \"16\u00d770=1120\"],\n  [\"18\u00d787=1566\", \"52\u00d739=2028\"],\n  [\"90\u00d769=6210\", \"48\u00d722=1056\"],\n  [\"49\u00d757=2793\", \"85\u00d797=8245\"],\n  [\"43\u00d739=1677\", \"21\u00d724=504\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the two-digit x two-digit multiplication answers in the table.\n# Each \"old\" value is unique in the document, so Find/Execute with\n# Replace:=wdReplaceAll against the whole document body is safe for every\n# pair even though a couple of the \"new\" values happen to collide with\n# \"old\" values used by other (already-processed) pairs.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"99\u00d722=2178\", \"97\u00d758=5626\"),\n    @(\"63\u00d754=3402\", \"33\u00d722=726\"),\n    @(\"21\u00d757=1197\", \"88\u00d753=4664\"),\n    @(\"26\u00d771=1846\", \"98\u00d796=9408\"),\n    @(\"39\u00d785=3315\", \"13\u00d789=1157\"),\n    @(\"48\u00d723=1104\", \"61\u00d713=793\"),\n    @(\"52\u00d746=2392\", \"14\u00d786=1204\"),\n    @(\"71\u00d777=5467\", \"82\u00d791=7462\"),\n    @(\"28\u00d727=756\", \"69\u00d734=2346\"),\n    @(\"65\u00d753=3445\", \"19\u00d786=1634\"),\n    @(\"38\u00d747=1786\", \"86\u00d792=7912\"),\n    @(\"73\u00d774=5402\", \"81\u00d773=5913\"),\n    @(\"32\u00d769=2208\", \"32\u00d759=1888\"),\n    @(\"36\u00d754=1944\", \"73\u00d756=4088\"),\n    @(\"15\u00d766=990\", \"98\u00d796=9408\"),\n    @(\"43\u00d731=1333\", \"38\u00d792=3496\"),\n    @(\"15\u00d738=570\", \"32\u00d769=2208\"),\n    @(\"87\u00d769=6003\", \"26\u00d732=832\"),\n    @(\"91\u00d772=6552\", \"87\u00d778=6786\"),\n    @(\"14\u00d716=224\", \"36\u00d726=936\"),\n    @(\"75\u00d776=5700\", \"16\u00d770=1120\"),\n    @(\"18\u00d787=1566\", \"52\u00d739=2028\"),\n    @(\"90\u00d769=6210\", \"48\u00d722=1056\"),\n    @(\"49\u00d757=2793\", \"85\u00d797=8245\"),\n    @(\"43\u00d739=1677\", \"21\u00d724=504\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
